# Update cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51.
# Leading apostrophe forces text entry (mirrors the original inline-string
# cells) so number-looking prices like "1.009" are not reinterpreted as numerics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.095.24"
$ws.Range("E2").Value = "  +2.02%  "

$ws.Range("D3").Value = "'1.950.02"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'328.38"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").Value = "'0.4853"
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").Value = "'0.4101"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.08242"
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("D10").Value = "'1.020"
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("D12").Value = "'1.948.64"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").Value = "'6.120"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").Value = "'7.350"
$ws.Range("E14").Value = "  +1.29%  "

$ws.Range("D15").Value = "'92.02"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "'0.06877"
$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").Value = "'0.00001042"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").Value = "'17.90"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("D20").Value = "'1.007"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "'30.065.19"
$ws.Range("E21").Value = "  +1.77%  "

$ws.Range("D22").Value = "'5.704"
$ws.Range("E22").Value = "  +1.18%  "

$ws.Range("D23").Value = "'11.99"
$ws.Range("E23").Value = "  +1.79%  "

$ws.Range("D24").Value = "'2.203"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("D25").Value = "'2.142.07"
$ws.Range("E25").Value = "  -2.05%  "

$ws.Range("D26").Value = "'6.597"
$ws.Range("E26").Value = "  -1.65%  "

$ws.Range("D27").Value = "'156.88"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "'20.15"
$ws.Range("E28").Value = "  +0.19%  "

$ws.Range("D29").Value = "'2.122"
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").Value = "'121.42"
$ws.Range("E30").Value = "  +0.59%  "

$ws.Range("D31").Value = "'1.026"
$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("D32").Value = "'0.09644"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").Value = "'5.649"
$ws.Range("E33").Value = "  +1.94%  "

$ws.Range("D34").Value = "'1.430"
$ws.Range("E34").Value = "  +2.66%  "

$ws.Range("D35").Value = "'3.556"
$ws.Range("E35").Value = "  -0.32%  "

$ws.Range("D36").Value = "'0.06554"
$ws.Range("E36").Value = "  +6.58%  "

$ws.Range("D37").Value = "'0.02308"
$ws.Range("E37").Value = "  +0.93%  "

$ws.Range("D38").Value = "'1.220"
$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("D39").Value = "'0.5981"
$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("D40").Value = "'10.78"
$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("D41").Value = "'7.988"
$ws.Range("E41").Value = "  -0.93%  "

$ws.Range("D42").Value = "'2.545"
$ws.Range("E42").Value = "  +5.49%  "

$ws.Range("D43").Value = "'0.1859"
$ws.Range("E43").Value = "  -0.43%  "

$ws.Range("D44").Value = "'12.52"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "'1.250"
$ws.Range("E45").Value = "  -2.48%  "

$ws.Range("D46").Value = "'0.07564"

$ws.Range("D47").Value = "'0.5592"
$ws.Range("E47").Value = "  -0.22%  "

$ws.Range("D48").Value = "'1.999"
$ws.Range("E48").Value = "  +1.93%  "

$ws.Range("D49").Value = "'118.21"
$ws.Range("E49").Value = "  +0.76%  "

$ws.Range("D50").Value = "'2.441"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("D51").Value = "'72.74"
$ws.Range("E51").Value = "  -0.16%  "
